$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 104
$ws.Range("I42").Value = 104
$ws.Range("K42").Value = 312
$ws.Range("M42").Value = -82

$ws.Range("H64").Value = 6245.4546
$ws.Range("J64").Value = 6325.5815
$ws.Range("L64").Value = 6325.5815
$ws.Range("N64").Value = -6821.5815

$ws.Range("H67").Value = 6245.4546
$ws.Range("J67").Value = 6325.5815
$ws.Range("L67").Value = 6325.5815
$ws.Range("N67").Value = -8041.5815

$ws.Range("H98").Value = 1025.3077
$ws.Range("I98").Value = 694.125
$ws.Range("K98").Value = 694.125
$ws.Range("M98").Value = 803.875

$ws.Range("H122").Value = 1025.3077
$ws.Range("I122").Value = 694.125
$ws.Range("K122").Value = 2082.375
$ws.Range("M122").Value = 367.625

$ws.Range("H135").Value = 1239.5
$ws.Range("I135").Value = 1047
$ws.Range("K135").Value = 9423
$ws.Range("M135").Value = -6888

$ws.Range("H138").Value = 6236.961
$ws.Range("I138").Value = 3649.6428
$ws.Range("K138").Value = 10948.9284
$ws.Range("M138").Value = -5808.928400000001

$ws.Range("H141").Value = 4044.7144
$ws.Range("I141").Value = 4226.8335
$ws.Range("J141").Value = 2952
$ws.Range("K141").Value = 12680.5005
$ws.Range("L141").Value = 8856
$ws.Range("M141").Value = -7500.500499999998
$ws.Range("N141").Value = -19216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 53552.633
$ws.Range("I2").Value = 72038.14
$ws.Range("K2").Value = 72038.14
$ws.Range("M2").Value = -71925.14

$ws.Range("H32").Value = 5051.8423
$ws.Range("I32").Value = 4837.0405
$ws.Range("K32").Value = 4837.0405
$ws.Range("M32").Value = -4550.0405

$ws.Range("H74").Value = 1566.6552
$ws.Range("I74").Value = 1690.1364
$ws.Range("J74").Value = 1178.5714
$ws.Range("K74").Value = 1690.1364
$ws.Range("L74").Value = 1178.5714
$ws.Range("M74").Value = -816.1364000000001
$ws.Range("N74").Value = -2926.5714

$ws.Range("H77").Value = 1566.6552
$ws.Range("I77").Value = 1690.1364
$ws.Range("J77").Value = 1178.5714
$ws.Range("K77").Value = 8450.682000000001
$ws.Range("L77").Value = 5892.857
$ws.Range("M77").Value = -4082.682000000001
$ws.Range("N77").Value = -14628.857

$ws.Range("H116").Value = 53552.633
$ws.Range("I116").Value = 72038.14
$ws.Range("K116").Value = 72038.14
$ws.Range("M116").Value = -69744.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 53552.633
$ws.Range("I3").Value = 72038.14
$ws.Range("K3").Value = 72038.14
$ws.Range("M3").Value = -71924.14

$ws.Range("H94").Value = 1009.59375
$ws.Range("I94").Value = 1135
$ws.Range("K94").Value = 1135
$ws.Range("M94").Value = -684

$ws.Range("H107").Value = 717320.1
$ws.Range("I107").Value = 2770
$ws.Range("J107").Value = 3337337.2
$ws.Range("K107").Value = 2770
$ws.Range("L107").Value = 3337337.2
$ws.Range("M107").Value = -850
$ws.Range("N107").Value = -3341177.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1432.6666
$ws.Range("I12").Value = 1299
$ws.Range("J12").Value = 1499.5
$ws.Range("K12").Value = 1299
$ws.Range("L12").Value = 1499.5
$ws.Range("M12").Value = -1129
$ws.Range("N12").Value = -1839.5

$ws.Range("H132").Value = 1708.6666
$ws.Range("I132").Value = 1384.1482
$ws.Range("J132").Value = 3169
$ws.Range("K132").Value = 4152.444600000001
$ws.Range("L132").Value = 9507
$ws.Range("M132").Value = -1622.444600000001
$ws.Range("N132").Value = -14567

$ws.Range("H134").Value = 239877.12
$ws.Range("I134").Value = 1813.4615
$ws.Range("K134").Value = 5440.3845
$ws.Range("M134").Value = -2905.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 499.25
$ws.Range("I12").Value = 498
$ws.Range("J12").Value = 499.66666
$ws.Range("K12").Value = 1494
$ws.Range("L12").Value = 1498.99998
$ws.Range("M12").Value = -1321
$ws.Range("N12").Value = -1844.99998

$ws.Range("H58").Value = 30333.223
$ws.Range("J58").Value = 31624.875
$ws.Range("L58").Value = 94874.625
$ws.Range("N58").Value = -95130.625

$ws.Range("H82").Value = 6307.1113
$ws.Range("I82").Value = 5999.5
$ws.Range("K82").Value = 17998.5
$ws.Range("M82").Value = -17592.5

$ws.Range("H85").Value = 6307.1113
$ws.Range("I85").Value = 5999.5
$ws.Range("K85").Value = 17998.5
$ws.Range("M85").Value = -16594.5

$ws.Range("H109").Value = 48980.816
$ws.Range("I109").Value = 1598.5
$ws.Range("K109").Value = 4795.5
$ws.Range("M109").Value = -3755.5

$ws.Range("H121").Value = 590070.4
$ws.Range("I121").Value = 1917.6
$ws.Range("J121").Value = 835134
$ws.Range("K121").Value = 5752.799999999999
$ws.Range("L121").Value = 2505402
$ws.Range("M121").Value = -4442.799999999999
$ws.Range("N121").Value = -2508022

$ws.Range("H122").Value = 1539.1666
$ws.Range("I122").Value = 884.7273
$ws.Range("J122").Value = 2092.923
$ws.Range("K122").Value = 7962.545700000001
$ws.Range("L122").Value = 18836.307
$ws.Range("M122").Value = -5512.545700000001
$ws.Range("N122").Value = -23736.307

$ws.Range("H137").Value = 2186.75
$ws.Range("I137").Value = 2264
$ws.Range("J137").Value = 1749
$ws.Range("K137").Value = 6792
$ws.Range("L137").Value = 5247
$ws.Range("M137").Value = -1692
$ws.Range("N137").Value = -15447

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 111129.8
$ws.Range("J75").Value = 111129.8
$ws.Range("L75").Value = 111129.8
$ws.Range("N75").Value = -112877.8

$ws.Range("H78").Value = 111129.8
$ws.Range("J78").Value = 111129.8
$ws.Range("L78").Value = 333389.4
$ws.Range("N78").Value = -342125.4

$ws.Range("H102").Value = 2326.6
$ws.Range("I102").Value = 1332.1072
$ws.Range("K102").Value = 1332.1072
$ws.Range("M102").Value = 289.8928000000001

$ws.Range("H126").Value = 3668.1667
$ws.Range("I126").Value = 2999
$ws.Range("J126").Value = 3802
$ws.Range("K126").Value = 8997
$ws.Range("L126").Value = 11406
$ws.Range("M126").Value = -6527
$ws.Range("N126").Value = -16346

$ws.Range("H132").Value = 90954.336
$ws.Range("I132").Value = 10690.5
$ws.Range("K132").Value = 32071.5
$ws.Range("M132").Value = -29541.5

$ws.Range("H134").Value = 49500
$ws.Range("J134").Value = 49500
$ws.Range("L134").Value = 148500
$ws.Range("N134").Value = -153570

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4604.0425
$ws.Range("I40").Value = 3846.3713
$ws.Range("K40").Value = 3846.3713
$ws.Range("M40").Value = -3710.3713

$ws.Range("H46").Value = 2361.818
$ws.Range("I46").Value = 2361
$ws.Range("K46").Value = 2361
$ws.Range("M46").Value = -2173

$ws.Range("H68").Value = 2379.6
$ws.Range("J68").Value = 2512
$ws.Range("L68").Value = 2512
$ws.Range("N68").Value = -4010

$ws.Range("H71").Value = 2379.6
$ws.Range("J71").Value = 2512
$ws.Range("L71").Value = 12560
$ws.Range("N71").Value = -20048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11153.75
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 11153.75
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H81").Value = 5981.8423
$ws.Range("I81").Value = 2775.6155
$ws.Range("J81").Value = 12928.667
$ws.Range("K81").Value = 5551.231
$ws.Range("L81").Value = 25857.334
$ws.Range("M81").Value = -4490.231
$ws.Range("N81").Value = -27979.334

$ws.Range("H84").Value = 5981.8423
$ws.Range("I84").Value = 2775.6155
$ws.Range("J84").Value = 12928.667
$ws.Range("K84").Value = 27756.155
$ws.Range("L84").Value = 129286.67
$ws.Range("M84").Value = -22452.155
$ws.Range("N84").Value = -139894.67

$ws.Range("H126").Value = 1455
$ws.Range("I126").Value = 1455
$ws.Range("K126").Value = 4365
$ws.Range("M126").Value = -1895

$ws.Range("H132").Value = 15983.25
$ws.Range("I132").Value = 2108.365
$ws.Range("J132").Value = 83223.08
$ws.Range("K132").Value = 6325.094999999999
$ws.Range("L132").Value = 249669.24
$ws.Range("M132").Value = -3795.094999999999
$ws.Range("N132").Value = -254729.24
